$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Cell values (entered in an order that matches the shared-strings layout of
# the target workbook: names down column A, then the header row, then the
# task descriptions column by column).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Thái Trung Đức"
$ws.Range("A3").Value = "Lê Chiêu Quốc"
$ws.Range("A4").Value = "Bùi Hữu Quý"
$ws.Range("A5").Value = "Phạm Như Ngọc Tuấn"

$ws.Range("B1").Value = "Tuần 1"
$ws.Range("C1").Value = "Tuần 2"

$ws.Range("D2").Value = "Thiết kế giao diện Mockup"
$ws.Range("D3").Value = "Vẽ System UseCase Diagram"
$ws.Range("D4").Value = "Tìm hiểu về đấu giá, đấu giá online(Khái niệm, quy trình, các chức năng cơ bản)"

$ws.Range("B2").Value = "Họp nhóm: Chọn đề tài, phân tích chức năng, viết Requirement Outline"
$ws.Range("E2").Value = "Viết tài liệu SRS"
$ws.Range("D5").Value = "Vẽ System Pages Flow"
$ws.Range("C2").Value = "Họp nhóm: Phân tích cụ thể từng chức năng, người sử dụng, luồng hoạt động của website => đặc tả chức năng"

# ---------------------------------------------------------------------------
# Merged cells
# ---------------------------------------------------------------------------
$ws.Range("B2:B5").Merge()
$ws.Range("E2:E3").Merge()
$ws.Range("C2:C5").Merge()
$ws.Range("C1:E1").Merge()

# ---------------------------------------------------------------------------
# Task-description cells (column D, regular font, vertical-center only)
# ---------------------------------------------------------------------------
$ws.Range("D2:D5").VerticalAlignment = -4108

# Header row (B1:E1): bold + centered both ways
$ws.Range("B1:E1").Font.Bold = $true
$ws.Range("B1:E1").HorizontalAlignment = -4108
$ws.Range("B1:E1").VerticalAlignment = -4108

# Names column (A2:A5): bold, vertical-center only
$ws.Range("A2:A5").Font.Bold = $true
$ws.Range("A2:A5").VerticalAlignment = -4108

# Meeting / milestone columns (B2:B5, C2:C5, E2:E3): centered both ways + wrap
$ws.Range("B2:B5").HorizontalAlignment = -4108
$ws.Range("B2:B5").VerticalAlignment = -4108
$ws.Range("B2:B5").WrapText = $true

$ws.Range("C2:C5").HorizontalAlignment = -4108
$ws.Range("C2:C5").VerticalAlignment = -4108
$ws.Range("C2:C5").WrapText = $true

$ws.Range("E2:E3").HorizontalAlignment = -4108
$ws.Range("E2:E3").VerticalAlignment = -4108
$ws.Range("E2:E3").WrapText = $true

# ---------------------------------------------------------------------------
# Column widths (closest achievable values to the original workbook's
# 17.625 / 20.125 / 29.75 / 61.5 / 12.5 character widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.832
$ws.Columns.Item(2).ColumnWidth = 19.332
$ws.Columns.Item(3).ColumnWidth = 29.003
$ws.Columns.Item(4).ColumnWidth = 60.665
$ws.Columns.Item(5).ColumnWidth = 11.665

# ---------------------------------------------------------------------------
# Sheet view / selection + page setup
# ---------------------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null
$ws.PageSetup.Orientation = 1

Write-Host "done"
